# Daily attendance processing - reorders the "Recorded By" (column G) list
# of contributors on every row where the token "System" appears alongside
# other recorders: the last recorder in the comma-separated list is moved
# to the front (a right-rotation of the list by one position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) { continue }

    $text = [string]$value
    if ($text.IndexOf(",") -lt 0) { continue }

    $parts = $text -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -lt 2) { continue }
    if (-not ($parts -contains "System")) { continue }

    $rotated = @($parts[$parts.Length - 1]) + $parts[0..($parts.Length - 2)]
    $newText = [string]::Join(", ", $rotated)

    $cell.Value = $newText
}
